$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# The workbook originally has 3 sheets: "Đơn sale chính", "Đơn phụ phẫu 1",
# and "Lương" (sheetId 3). We need to:
#   1. Duplicate the existing "Lương" sheet, placing the copy right after
#      it, so the copy becomes the new sheet (sheetId 4) that keeps the
#      name "Lương" (with an inserted bonus row + updated totals).
#   2. Turn the original third sheet into the new "Thưởng" sheet, with a
#      fresh bonus/penalty table.
# ------------------------------------------------------------------

$luong = $wb.Worksheets.Item(3)

# Step 1: duplicate "Lương" -> placed immediately after itself.
$luong.Copy($null, $luong)
$luongNew = $wb.Worksheets.Item(4)
$luongNew.Name = "Luong_tmp"

# Step 2: insert the new "Thưởng tại CẦN THƠ" row at row 11 of the copy,
# shifting every following row down by one (34 rows -> 35 rows).
$luongNew.Rows.Item(11).Insert()
$luongNew.Range("A11").Value = "Thưởng tại CẦN THƠ"
$luongNew.Range("B11").Value = 4000000

# Step 3: update the two grand-total rows, which are now at row 32
# ("Tổng lương tại CẦN THƠ") and row 35 ("Tổng lương tại HỆ THỐNG"),
# by adding in the new bonus amount.
$tongCanTho = $luongNew.Range("B32").Value2
$luongNew.Range("B32").Value = $tongCanTho + 4000000

$tongHeThong = $luongNew.Range("B35").Value2
$luongNew.Range("B35").Value = $tongHeThong + 4000000

# ------------------------------------------------------------------
# Step 4: turn the original "Lương" sheet into the "Thưởng" sheet: wipe
# its contents and write the bonus/penalty report table.
# ------------------------------------------------------------------
$luong.Cells.Clear()

$headers = @("Ngày phát sinh", "notion id", "Tiền tố", "Mã thưởng phạt", "id nhân sự", "Cơ sở", "Loại", "Lượng thưởng phạt", "Lí do", "Họ và tên")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $col = $i + 1
    $luong.Cells.Item(1, $col).Value = $headers[$i]
}

# Force column A to text so the date-looking string isn't reinterpreted
# as a real date serial number.
$luong.Range("A2").NumberFormat = "@"

$luong.Range("A2").Value = "08-05-2024"
$luong.Range("B2").Value = "f1df828e-4b4e-4cec-93fb-1de41e0d82f5"
$luong.Range("C2").Value = "TP"
$luong.Range("D2").Value = 9
$luong.Range("E2").Value = "e49d0ce3-124d-4e4b-b377-be2139cde3f5"
$luong.Range("F2").Value = "CẦN THƠ"
$luong.Range("G2").Value = "Thưởng"
$luong.Range("H2").Value = 4000000
$luong.Range("J2").Value = "Lâm Hoàng Phú"

$luong.Range("C3").Value = "Tổng"
$luong.Range("D3").Value = 1
$luong.Range("H3").Value = 4000000

# Now rename the sheets into their final names.
$luong.Name = "Thưởng"
$luongNew.Name = "Lương"
